# "quitar el componente de imprimir pantalla"
# Close out two "hallazgos" (findings) rows that belonged to the DDC
# print-screen feature: mark the finding as resolved (Cerrada / APROBADO),
# correct the status of another finding, normalize the alignment of the
# description column for that block, and re-point the active sheet/
# selection at the "hallazgos" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hallazgos")

# H7: estado cambia de "REVISAR" a "En proceso"
$ws.Cells.Item(7, 8).Value = "En proceso"

# Fila 14: se cierra el hallazgo y queda aprobado
$ws.Cells.Item(14, 7).Value = "Cerrada"
$ws.Cells.Item(14, 8).Value = "APROBADO"

# Fila 15 ("Quitar Impresión de matriz (DDC)"): se cierra el hallazgo y queda aprobado
$ws.Cells.Item(15, 7).Value = "Cerrada"
$ws.Cells.Item(15, 8).Value = "APROBADO"

# Centrar la columna D (Descripción del error) para las filas 12 a 19
$colD = $ws.Range("D12:D19")
$colD.HorizontalAlignment = -4108
$colD.VerticalAlignment = -4108

# Alinear el formato de la fila 20 con el resto de la tabla (tomando la
# fila 19 como referencia) y luego centrar su celda D
$srcRow = $ws.Range("B19:G19")
$dstRow = $ws.Range("B20:G20")
$srcRow.Copy()
$dstRow.PasteSpecial(-4122)
$excel.CutCopyMode = $false

$d20 = $ws.Cells.Item(20, 4)
$d20.HorizontalAlignment = -4108
$d20.VerticalAlignment = -4108

# Se deja activa la hoja "hallazgos" (en vez de "RF") con H16 seleccionada
$ws.Activate()
[void]$ws.Range("H16").Select()
